$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# RGB color matching the "D9E2F3" (accent1, 33% tint) shading already used
# elsewhere in this table for "section separator" rows.
$shadingColor = 0xD9 + (0xE2 * 256) + (0xF3 * 65536)

function Fill-Row($row, $texts) {
    for ($i = 1; $i -le $row.Cells.Count; $i++) {
        if ($texts[$i - 1] -ne $null) {
            $row.Cells.Item($i).Range.Text = $texts[$i - 1]
        }
    }
}

function Shade-Row($row) {
    for ($i = 1; $i -le $row.Cells.Count; $i++) {
        $row.Cells.Item($i).Shading.BackgroundPatternColor = $shadingColor
    }
}

# New row: "get Admin/allJobItems/:id/:token" -> jobItems -> "Returns map of
# all existing jobitems"
$row1 = $t.Rows.Add()
Fill-Row $row1 @("get", "Admin/allJobItems/:id/:token", $null, "jobItems", "Returns map of all existing jobitems")
Shade-Row $row1

# New row: "get Admin/allCompanies/:id/:token" -> Companies -> "Returns map
# of all existing companies"
$row2 = $t.Rows.Add()
Fill-Row $row2 @("get", "Admin/allCompanies/:id/:token", $null, "Companies", "Returns map of all existing companies")
Shade-Row $row2

Write-Host "Table now has" $t.Rows.Count "rows"
